$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells (I1, J1), matching the existing header style (H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new data columns I (I0) and J (IF) for rows 2-15
$data = @{
    2  = @(1,4)
    3  = @(1,6)
    4  = @(1,5)
    5  = @(1,5)
    6  = @(1,6)
    7  = @(1,6)
    8  = @(1,5)
    9  = @(1,5)
    10 = @(4,4)
    11 = @(3,7)
    12 = @(1,4)
    13 = @(1,4)
    14 = @(1,4)
    15 = @(8,8)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
